# Reports commited for 3rd June
# - Test Steps (Creator_platform rows TS029..TS064): mark Result1 as PASS
# - Creator_platform: mark Result as PASS, roll test account emails forward
# - Reset cell selections left over from the authoring session

$wb = $excel.ActiveWorkbook

$wsCases     = $wb.Worksheets.Item("Test Cases")
$wsSteps     = $wb.Worksheets.Item("Test Steps")
$wsCreator   = $wb.Worksheets.Item("Creator_platform")

# ---------------------------------------------------------------------
# Test Steps: rows that previously held long Selenium FAIL stack traces
# now simply read PASS (rows 31, 44 and 49 were already PASS).
# ---------------------------------------------------------------------
$passRows = @(30,32,33,34,35,36,37,38,39,40,41,42,43,45,46,47,48,50,51,52,53,54,55,56,57,58,59,60,61,62,63,64,65)
foreach ($r in $passRows) {
    $wsSteps.Range("H$r").Value = "PASS"
}
$wsSteps.Columns("H").AutoFit()

# ---------------------------------------------------------------------
# Creator_platform: new pair of test accounts for this run, and the
# single result cell on the sheet now reads PASS.
# ---------------------------------------------------------------------
$wsCreator.Range("A2").Value = "'genvideotest1247@gmail.com"
$wsCreator.Range("C2").Value = "'genvideotest1248@gmail.com"
$wsCreator.Range("R2").Value = "PASS"
$wsCreator.Columns("R").AutoFit()

# ---------------------------------------------------------------------
# Leftover cell selections from the authoring session
# ---------------------------------------------------------------------
$wsCases.Range("C2").Select()
$wsSteps.Range("E1").Select()

$wsCreator.Activate()
$wsCreator.Range("A3").Select()
